$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.457.92"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "2.302.06"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.965"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "2.651.90"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "2.291.94"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").Value = "42.423.30"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("E19").Value = "  -2.18%  "

$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "276.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +20.90%  "

$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("E28").Value = "  +3.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0873"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("E34").Value = "  +3.96%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.06%  "

$ws.Range("E37").Value = "  +3.42%  "

$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("E41").Value = "  +2.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.76%  "

$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "81.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.18%  "

$ws.Range("E47").Value = "  -2.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.17%  "

$ws.Range("D51").Value = "1.590.12"
$ws.Range("E51").Value = "  +1.26%  "
